# "Grid Search with Keras"
#
# 1) Slide 11 ("Scikit-Learn Wrapper for Keras") - extend the bold heading
#    in the body textbox so it reads "... Cross-validation of a Deep Net".
# 2) The cached/auto-updating footer date fields (type="datetime1" /
#    "datetimeFigureOut") on the slide master and every slide layout that
#    carries a date placeholder advanced by one day: 1/29/22 -> 1/30/22.

$p = $ppt.ActivePresentation

# --- 1. Update the heading text on slide 11 -------------------------------
$slide = $p.Slides.Item(11)
$titleShape = $slide.Shapes.Item(3)
$titleRange = $titleShape.TextFrame.TextRange
$firstPara = $titleRange.Paragraphs(1, 1)
$firstPara.Characters(1, $firstPara.Length).Text = "Facilitate Hyperparameter Tuning and Cross-validation of a Deep Net"

# --- 2. Roll the footer date forward on the slide master and layouts ------
$newDate = "1/30/22"
$oldDate = "1/29/22"

$master = $p.SlideMaster

$masterDateShape = $master.Shapes.Item("Date Placeholder 3")
if ($masterDateShape.TextFrame.TextRange.Text -eq $oldDate) {
    $masterDateShape.TextFrame.TextRange.Text = $newDate
}

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shp = $layout.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}
